$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column C (the "Förändrad" column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

# Update every cell in column C from row 2 to the last row that currently
# holds the old date serial (45190) to the new date serial (45192).
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45192
